# Atualizado por script em 02-12-2023 02:45
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap the data of rows 128 and 129 (columns F..V); Indice (A) and data_partida (E) stay put ---
$cols = @("F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V")

$row128 = @()
$row129 = @()
foreach ($c in $cols) {
    $row128 += ,($ws.Range($c + "128").Value2)
    $row129 += ,($ws.Range($c + "129").Value2)
}

for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "128").Value = $row129[$i]
    $ws.Range($cols[$i] + "129").Value = $row128[$i]
}

# --- Append new rows 131-133, copying the formatting (style) of the A and E columns from row 130 ---
$ws.Range("A130").Copy()
$ws.Range("A131:A133").PasteSpecial(-4122)

$ws.Range("E130").Copy()
$ws.Range("E131:E133").PasteSpecial(-4122)

# The "temporada" column (D) holds a numeric-looking string ("2023"); assigning
# it through .Value would be auto-coerced to a number, so instead copy the
# already-text-typed D130 cell's value (xlPasteValues) into each new row.
$ws.Range("D130").Copy()
$ws.Range("D131").PasteSpecial(-4163)
$ws.Range("D130").Copy()
$ws.Range("D132").PasteSpecial(-4163)
$ws.Range("D130").Copy()
$ws.Range("D133").PasteSpecial(-4163)

# Row 131
$ws.Range("A131").Value = 130
$ws.Range("B131").Value = "paraguay"
$ws.Range("C131").Value = "primera-division"
$ws.Range("E131").Value = 45261.97916666666
$ws.Range("F131").Value = "Cerro Porteno"
$ws.Range("G131").Value = 4
$ws.Range("H131").Value = "Guarani"
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 1.63
$ws.Range("K131").Value = "25/11/2023 23:42"
$ws.Range("L131").Value = 1.85
$ws.Range("M131").Value = "01/12/2023 23:28"
$ws.Range("N131").Value = 4.04
$ws.Range("O131").Value = "25/11/2023 23:42"
$ws.Range("P131").Value = 3.84
$ws.Range("Q131").Value = "01/12/2023 23:28"
$ws.Range("R131").Value = 5.39
$ws.Range("S131").Value = "25/11/2023 23:42"
$ws.Range("T131").Value = 4.04
$ws.Range("U131").Value = "01/12/2023 23:28"
$ws.Range("V131").Value = "https://www.betexplorer.com/football/paraguay/primera-division/cerro-porteno-guarani/Ag2Jqswg/"

# Row 132
$ws.Range("A132").Value = 131
$ws.Range("B132").Value = "paraguay"
$ws.Range("C132").Value = "primera-division"
$ws.Range("E132").Value = 45261.97916666666
$ws.Range("F132").Value = "General Caballero JLM"
$ws.Range("G132").Value = 0
$ws.Range("H132").Value = "Olimpia Asuncion"
$ws.Range("I132").Value = 1
$ws.Range("J132").Value = 4
$ws.Range("K132").Value = "25/11/2023 23:42"
$ws.Range("L132").Value = 3.59
$ws.Range("M132").Value = "01/12/2023 23:27"
$ws.Range("N132").Value = 3.45
$ws.Range("O132").Value = "25/11/2023 23:42"
$ws.Range("P132").Value = 3.12
$ws.Range("Q132").Value = "01/12/2023 23:27"
$ws.Range("R132").Value = 1.91
$ws.Range("S132").Value = "25/11/2023 23:42"
$ws.Range("T132").Value = 2.24
$ws.Range("U132").Value = "01/12/2023 23:27"
$ws.Range("V132").Value = "https://www.betexplorer.com/football/paraguay/primera-division/general-caballero-jlm-olimpia-asuncion/x2os3VwE/"

# Row 133
$ws.Range("A133").Value = 132
$ws.Range("B133").Value = "paraguay"
$ws.Range("C133").Value = "primera-division"
$ws.Range("E133").Value = 45261.97916666666
$ws.Range("F133").Value = "Sp. Luqueno"
$ws.Range("G133").Value = 1
$ws.Range("H133").Value = "Nacional Asuncion"
$ws.Range("I133").Value = 1
$ws.Range("J133").Value = 2.65
$ws.Range("K133").Value = "25/11/2023 23:42"
$ws.Range("L133").Value = 3.1
$ws.Range("M133").Value = "01/12/2023 23:26"
$ws.Range("N133").Value = 3.12
$ws.Range("O133").Value = "25/11/2023 23:42"
$ws.Range("P133").Value = 3
$ws.Range("Q133").Value = "01/12/2023 23:27"
$ws.Range("R133").Value = 2.89
$ws.Range("S133").Value = "25/11/2023 23:42"
$ws.Range("T133").Value = 2.58
$ws.Range("U133").Value = "01/12/2023 23:26"
$ws.Range("V133").Value = "https://www.betexplorer.com/football/paraguay/primera-division/sp-luqueno-nacional-asuncion/W23FpNNn/"
